# Add the new "PhysicalServerModelName" property column.
#
# The "VMs" sheet currently has headers:
#   A1 Physical Hosts | B1 Hostname | C1 VirtualizationNode | D1 CappedCPU
#
# We insert a new column before D for "PhysicalServerModelName", which
# pushes the existing "CappedCPU" column from D to E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before D - shifts "CappedCPU" (old D) to E.
$ws.Columns("D:D").Insert()

# Populate the header of the newly inserted column.
$ws.Range("D1").Value = "PhysicalServerModelName"

# Widen the new "PhysicalServerModelName" column to fit its longer header,
# and restore "CappedCPU"'s original column width now that it lives in E.
$ws.Columns("D:D").ColumnWidth = 22.25
$ws.Columns("E:E").ColumnWidth = 11.084

# Reset the active selection back to A1.
$ws.Range("A1").Select()
